$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Current state (before this edit) ---
# Row 4: torneo | Prejuveniles | caballeros | 1 | Luján Martínez, Benjamín | 75 | 84 | 159
# Row 5: torneo | Prejuveniles | caballeros | 2 | Petric, Juan Cruz        | 81 | 79 | 160
#
# --- Target state (after this edit) ---
# Row 4: torneo | Prejuveniles | caballeros | 1 | Kern Pascuali, Juan Daniel | 92 | (blank) | 92
# Row 5: torneo | Prejuveniles | caballeros | 2 | Luján Martínez, Benjamín   | 75 | 84      | 159
# Row 6: torneo | Prejuveniles | caballeros | 3 | Petric, Juan Cruz          | 81 | 79      | 160
#
# i.e. a new result (Kern Pascuali, Juan Daniel) is inserted at the top of the
# "Prejuveniles" group, bumping the existing Luján Martínez row down to
# position 2 (row 5) and the existing Petric row down to position 3 (row 6).

# Remember row 4's current contents - they get carried down to the new row 5.
$carryA = $ws.Range("A4").Value2
$carryB = $ws.Range("B4").Value2
$carryC = $ws.Range("C4").Value2
$carryE = $ws.Range("E4").Value2
$carryF = $ws.Range("F4").Value2
$carryG = $ws.Range("G4").Value2
$carryH = $ws.Range("H4").Value2

# Insert a new row above the current row 5, pushing the Petric row down to row 6.
$ws.Rows("5:5").Insert()

# The Petric row (now row 6) moves from position 2 to position 3.
$ws.Range("D6").Value2 = 3

# The newly inserted row 5 gets what used to be in row 4 (Luján Martínez), still
# at position 2.
$ws.Range("A5").Value2 = $carryA
$ws.Range("B5").Value2 = $carryB
$ws.Range("C5").Value2 = $carryC
$ws.Range("D5").Value2 = 2
$ws.Range("E5").Value2 = $carryE
$ws.Range("F5").Value2 = $carryF
$ws.Range("G5").Value2 = $carryG
$ws.Range("H5").Value2 = $carryH

# Row 4 itself now holds the new player's result (position 1, day 2 blank).
$ws.Range("E4").Value2 = "Kern Pascuali, Juan Daniel"
$ws.Range("F4").Value2 = 92
$ws.Range("G4").Value2 = ""
$ws.Range("H4").Value2 = 92
